# "added concept part to documentation"
#
# The Zeitaufzeichnung (time-log) sheet gets one new entry appended just
# above the "total" row: a "concept part" work session on 30.12.2024
# worth 9 hours. That pushes the running SUM() total from 233 to 242
# (it recalculates automatically).
#
# Previously row 54 ("29.12.2024" / 10h / "Kommentieren von Code, ...")
# was the last data row, so it carried the sheet's "last row" emphasis
# formatting; now that the new row 55 is the last data row, that emphasis
# moves down to row 55 and row 54 reverts to the regular/default look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54 is no longer the last entry - drop its special emphasis back to
# the regular (unstyled) look used by every other data row.
$ws.Range("A54").Style = "Normal"
$ws.Range("C54").Style = "Normal"

# New last entry: concept part of the documentation, 30.12.2024, 9 hours.
$ws.Range("A55").Value = "30.12.2024"
$ws.Range("B55").Value = 9
